# Updated temperature (column K: average_county_temperature) with NOAA data,
# which also changes the downstream worst_ashp_cop / best_ashp_cop columns
# (R and S) that depend on it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mapping of row -> new values for columns K (temperature), R (worst_ashp_cop),
# S (best_ashp_cop). Rows without R/S entries only had K recomputed.
$updates = @(
    @{ Row = 2;  K = 13.75752314814816; R = 1.722630989917367; S = 1.865269081797952 },
    @{ Row = 3;  K = 13.75752314814816; R = 1.803186500133452; S = 1.964569140204562 },
    @{ Row = 7;  K = 13.75752314814816; R = 1.722630989917367; S = 1.865269081797952 },
    @{ Row = 8;  K = 13.75752314814816; R = 1.803186500133452; S = 1.964569140204562 },
    @{ Row = 9;  K = 3.38888888888889;  R = 1.578134831460674; S = 1.695036674816626 },
    @{ Row = 10; K = 3.38888888888889;  R = 1.642853739876131; S = 1.773501823866597 },
    @{ Row = 13; K = 12.93898809523811 },
    @{ Row = 14; K = 12.93898809523811 },
    @{ Row = 15; K = 12.93898809523811 },
    @{ Row = 18; K = 19.79629629629628; R = 1.911855479578636; S = 2.09608909874769 },
    @{ Row = 19; K = 19.79629629629628 },
    @{ Row = 20; K = 19.79629629629628; R = 1.819666609086197; S = 1.981148790245761 },
    @{ Row = 21; K = 13.75752314814816; R = 1.722630989917367; S = 1.865269081797952 },
    @{ Row = 22; K = 13.75752314814816; R = 1.803186500133452; S = 1.964569140204562 },
    @{ Row = 23; K = 3.38888888888889;  R = 1.578134831460674; S = 1.695036674816626 },
    @{ Row = 24; K = 3.38888888888889;  R = 1.642853739876131; S = 1.773501823866597 },
    @{ Row = 25; K = 3.38888888888889 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("K$r").Value = $u.K
    if ($u.ContainsKey("R")) {
        $ws.Range("R$r").Value = $u.R
    }
    if ($u.ContainsKey("S")) {
        $ws.Range("S$r").Value = $u.S
    }
}
